# Zeiterfassung.xlsx - add a new weekly time-tracking block (03.11.17)
# mirroring the existing "02.11.17" block (rows 41-47), appended after
# the last block (rows 33-47) as rows 48-55.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the previous full block (separator row 40 + data rows 41-47,
# 8 rows total) straight down to rows 48-55. Range.Copy(Destination)
# carries over both values and cell formatting (styles) in one shot.
$ws.Range("A40:D47").Copy($ws.Range("A48:D55"))

# Row 49 is the new block's header row: fix the date to the new week.
$ws.Range("A49").Value = "03.11.17"

# Row 53 ("- Logik"): Finn (column C) logged "1".
$ws.Range("C53").Value = "1"

# Row 55 ("Recherche" totals row): Finn (column C) totals 1 (numeric).
$ws.Range("C55").Value = 1
